$d = $word.ActiveDocument

function Set-RangeXml($range, $innerRunsXml) {
    # Re-seat the range on a fresh Document.Range object. Range objects
    # handed back from Paragraph.Range (even after MoveEnd) make
    # Range.InsertXML() *append* its payload instead of replacing the
    # range's contents; rebuilding the same [start,end) span via
    # Document.Range() makes InsertXML replace correctly.
    $fresh = $d.Range($range.Start, $range.End)
    $openXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
               '<pkg:part pkg:name="/word/document.xml">' + `
               '<pkg:xmlData>' + `
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
               '<w:body><w:p>' + $innerRunsXml + '</w:p></w:body>' + `
               '</w:document>' + `
               '</pkg:xmlData></pkg:part></pkg:package>'
    $fresh.InsertXML($openXml)
}

# 1. Update the timestamp in the Date-styled paragraph near the top of the
#    document: "10:33:58 AM" -> "08:23:56 PM" (the date portion is unchanged).
$dateTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*June  23, 2021*") {
        $dateTarget = $p
    }
}

if ($dateTarget -ne $null) {
    $r = $dateTarget.Range
    $r.MoveEnd(1, -1) | Out-Null
    $newDateRun = '<w:r><w:t xml:space="preserve">June  23, 2021 (08:23:56 PM)</w:t></w:r>'
    Set-RangeXml $r $newDateRun
}

# 2. Collapse the syntax-highlighted "Cannot implicitly convert type ..."
#    source-code line (originally ten runs, each with its own Pygments
#    token rStyle) into a single run styled as plain "Verbatim Char".
$codeTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Cannot implicitly convert type*") {
        $codeTarget = $p
    }
}

if ($codeTarget -ne $null) {
    $r2 = $codeTarget.Range
    $r2.MoveEnd(1, -1) | Out-Null
    $mergedText = "Cannot implicitly convert type 'float' to 'int'. An explicit conversion exists (are you missing a cast?)"
    $newCodeRun = '<w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr>' + `
                  '<w:t xml:space="preserve">' + $mergedText + '</w:t></w:r>'
    Set-RangeXml $r2 $newCodeRun
}
